$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.199.78"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = "'3.675.94"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.30%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'595.87"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = "'165.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.06%  '
$ws.Range('D7').Value = "'3.674.46"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.20%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'0.531"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').Value = "'0.165"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.14%  '
$ws.Range('D11').Value = "'6.24"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = "'0.457"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = "'37.89"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').Value = "'0.0000244"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').Value = "'4.296.06"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').Value = "'3.687.54"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.25%  '
$ws.Range('D17').Value = "'68.157.87"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = "'7.22"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').Value = "'0.114"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = "'17.09"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.34%  '
$ws.Range('D21').Value = "'488.84"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'9.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').Value = "'0.718"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').Value = "'84.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').Value = "'0.0000141"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').Value = "'2.28"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.98%  '
$ws.Range('D27').Value = "'12.15"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').Value = "'10.02"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = "'2.90"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('B31').Value = "'ImmutableX"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'2.37"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.05%  '
$ws.Range('B32').Value = "'NEARProtocol"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'7.81"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = "'31.19"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.67%  '
$ws.Range('D34').Value = "'3.822.79"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').Value = "'3.623.40"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = "'0.994"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.43%  '
$ws.Range('D39').Value = "'5.73"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('D40').Value = "'0.131"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').Value = "'0.320"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').Value = "'431.70"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('D43').Value = "'48.47"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('D44').Value = "'1.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').Value = "'2.81"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('D46').Value = "'8.34"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.90%  '
$ws.Range('D48').Value = "'40.30"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('D49').Value = "'141.43"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').Value = "'2.738.06"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.49%  '
$ws.Range('D51').Value = "'0.0348"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.01%  '
